# Auto-applied updates to Leve profit sheets (market price refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 112.666664
$ws.Range("I4").Value = 131.66667
$ws.Range("J4").Value = 74.666664
$ws.Range("K4").Value = 131.66667
$ws.Range("L4").Value = 74.666664
$ws.Range("M4").Value = -17.66667000000001
$ws.Range("N4").Value = -302.666664
$ws.Range("H18").Value = 693.5
$ws.Range("J18").Value = 725
$ws.Range("L18").Value = 725
$ws.Range("N18").Value = -1293
$ws.Range("H52").Value = 2121
$ws.Range("J52").Value = 2121
$ws.Range("L52").Value = 6363
$ws.Range("N52").Value = -6683
$ws.Range("H100").Value = 4000
$ws.Range("I100").Value = 4000
$ws.Range("K100").Value = 4000
$ws.Range("M100").Value = -3459
$ws.Range("H106").Value = 2554.2727
$ws.Range("I106").Value = 2639.7
$ws.Range("K106").Value = 2639.7
$ws.Range("M106").Value = -2008.7
$ws.Range("H113").Value = 21900
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2219.8
$ws.Range("I45").Value = 1100
$ws.Range("K45").Value = 1100
$ws.Range("M45").Value = -723
$ws.Range("H61").Value = 1638.4667
$ws.Range("I61").Value = 1698.6154
$ws.Range("K61").Value = 1698.6154
$ws.Range("M61").Value = -1486.6154
$ws.Range("H122").Value = 889
$ws.Range("I122").Value = 889
$ws.Range("K122").Value = 2667
$ws.Range("M122").Value = -217
$ws.Range("H132").Value = 1274.6451
$ws.Range("I132").Value = 1195
$ws.Range("K132").Value = 3585
$ws.Range("M132").Value = -1055
$ws.Range("H136").Value = 1638.4667
$ws.Range("I136").Value = 1698.6154
$ws.Range("K136").Value = 5095.8462
$ws.Range("M136").Value = -2545.8462

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1553.3077
$ws.Range("I86").Value = 1309
$ws.Range("J86").Value = 2897
$ws.Range("K86").Value = 1309
$ws.Range("L86").Value = 2897
$ws.Range("M86").Value = -186
$ws.Range("N86").Value = -5143
$ws.Range("H89").Value = 1553.3077
$ws.Range("I89").Value = 1309
$ws.Range("J89").Value = 2897
$ws.Range("K89").Value = 6545
$ws.Range("L89").Value = 14485
$ws.Range("M89").Value = -929
$ws.Range("N89").Value = -25717
$ws.Range("H94").Value = 2985.6667
$ws.Range("I94").Value = 2985.6667
$ws.Range("K94").Value = 2985.6667
$ws.Range("M94").Value = -2534.6667
$ws.Range("H105").Value = 3581
$ws.Range("I105").Value = 3533.7
$ws.Range("J105").Value = 3699.25
$ws.Range("K105").Value = 3533.7
$ws.Range("L105").Value = 3699.25
$ws.Range("M105").Value = -1786.7
$ws.Range("N105").Value = -7193.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1343.1
$ws.Range("I16").Value = 1353.5
$ws.Range("K16").Value = 1353.5
$ws.Range("M16").Value = -1066.5
$ws.Range("H31").Value = 2247.5405
$ws.Range("I31").Value = 1700.3334
$ws.Range("K31").Value = 1700.3334
$ws.Range("M31").Value = -1405.3334
$ws.Range("H34").Value = 2247.5405
$ws.Range("I34").Value = 1700.3334
$ws.Range("K34").Value = 1700.3334
$ws.Range("M34").Value = -1498.3334
$ws.Range("H113").Value = 1343.1
$ws.Range("I113").Value = 1353.5
$ws.Range("K113").Value = 1353.5
$ws.Range("M113").Value = 816.5
$ws.Range("H132").Value = 1730.4445
$ws.Range("I132").Value = 1595.6666
$ws.Range("K132").Value = 4786.9998
$ws.Range("M132").Value = -2256.9998
$ws.Range("H134").Value = 2929.125
$ws.Range("I134").Value = 1880.5
$ws.Range("K134").Value = 5641.5
$ws.Range("M134").Value = -3106.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2566.7144
$ws.Range("I80").Value = 2396
$ws.Range("J80").Value = 2993.5
$ws.Range("K80").Value = 2396
$ws.Range("L80").Value = 2993.5
$ws.Range("M80").Value = -1398
$ws.Range("N80").Value = -4989.5
$ws.Range("H83").Value = 2566.7144
$ws.Range("I83").Value = 2396
$ws.Range("J83").Value = 2993.5
$ws.Range("K83").Value = 11980
$ws.Range("L83").Value = 14967.5
$ws.Range("M83").Value = -6988
$ws.Range("N83").Value = -24951.5
$ws.Range("H102").Value = 3044.25
$ws.Range("I102").Value = 3044.25
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3044.25
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1422.25
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 2466.5625
$ws.Range("J122").Value = 6499
$ws.Range("L122").Value = 19497
$ws.Range("N122").Value = -24397
$ws.Range("H124").Value = 56975
$ws.Range("J124").Value = 56975
$ws.Range("L124").Value = 56975
$ws.Range("N124").Value = -66795

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 5000
$ws.Range("K5").Value = 5000
$ws.Range("M5").Value = -4887
$ws.Range("H7").Value = 68994.5
$ws.Range("I7").Value = 68994.5
$ws.Range("K7").Value = 68994.5
$ws.Range("M7").Value = -68882.5
$ws.Range("H13").Value = 2006
$ws.Range("I13").Value = 2006
$ws.Range("K13").Value = 2006
$ws.Range("M13").Value = -1866
$ws.Range("H40").Value = 2997
$ws.Range("I40").Value = 2997
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2997
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2861
$ws.Range("N40").ClearContents()
$ws.Range("H64").Value = 16000
$ws.Range("I64").Value = 12000
$ws.Range("J64").Value = 20000
$ws.Range("K64").Value = 12000
$ws.Range("L64").Value = 20000
$ws.Range("M64").Value = -11775
$ws.Range("N64").Value = -20450
$ws.Range("H67").Value = 16000
$ws.Range("I67").Value = 12000
$ws.Range("J67").Value = 20000
$ws.Range("K67").Value = 12000
$ws.Range("L67").Value = 20000
$ws.Range("M67").Value = -11220
$ws.Range("N67").Value = -21560
$ws.Range("H126").Value = 68994.5
$ws.Range("I126").Value = 68994.5
$ws.Range("K126").Value = 206983.5
$ws.Range("M126").Value = -204513.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 19375
$ws.Range("H66").Value = 19375
$ws.Range("H126").Value = 1249.75
$ws.Range("I126").Value = 999.5
$ws.Range("K126").Value = 2998.5
$ws.Range("M126").Value = -528.5
